$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the payment-code references from DISK230100005 to DISK230300012 ---
# F2 holds the multi-line "PREPARATION" text containing the old code.
$nl = [char]10
$ws.Range("F2").Value = "Username : 37841;" + $nl + "Password : bni1234;" + $nl + "Role : RL09 - Penyelia Settlement;" + $nl + "Kode Pembayaran : DISK230300012"

# N2 holds the standalone "KODE_PEMBAYARAN" value.
$ws.Range("N2").Value = "DISK230300012"

# --- Alignment touch-ups on row 2 ---
# B2, D2, E2 switch from centered to left-aligned (matching F2's existing left alignment).
$ws.Range("B2").HorizontalAlignment = -4131   # xlLeft
$ws.Range("D2").HorizontalAlignment = -4131   # xlLeft
$ws.Range("E2").HorizontalAlignment = -4131   # xlLeft

# F2 gains an explicit left horizontal alignment (previously unset/general).
$ws.Range("F2").HorizontalAlignment = -4131   # xlLeft

# --- Update the active selection / scroll position to reflect the new view ---
$ws.Range("O2").Select()
